# fix update perhitungan btn indonesia
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the figures in rows 2-4 (penghasilan / pengeluaran columns)
$ws.Range("A2").Value = 105000000
$ws.Range("B2").Value = 302000

$ws.Range("A3").Value = 500400000
$ws.Range("B3").Value = 3200000

$ws.Range("A4").Value = 3004000000
$ws.Range("B4").Value = 5004000

# Remove the now-unused trailing rows (5-8) so the sheet only covers A1:C4
$ws.Range("A5:C8").EntireRow.Delete()

# Restore the outline summary levels recorded on the sheet (row/col) to match
# the shrunk data range, without leaving any stray per-row/per-col grouping.
$ws.Outline.ShowLevels(3, 2)

# Match the author's final selection
$ws.Range("B14").Select()
